# Updated cryptos list on Tue Feb  6 21:48:39 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.240.90"
$ws.Range("E2").Value = "  +1.72%  "

# Row 3
$ws.Range("D3").Value = "2.384.94"
$ws.Range("E3").Value = "  +4.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.81%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.09"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.67%  "

# Row 7
$ws.Range("E7").Value = "  +0.29%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.98%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.23"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.52%  "

# Row 11
$ws.Range("E11").Value = "  +1.18%  "

# Row 12
$ws.Range("E12").Value = "  +2.46%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.44"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.02%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.80"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.04%  "

# Row 15
$ws.Range("D15").Value = "2.754.24"
$ws.Range("E15").Value = "  +4.14%  "

# Row 16
$ws.Range("D16").Value = "2.373.52"
$ws.Range("E16").Value = "  +3.52%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.810"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.92%  "

# Row 18
$ws.Range("D18").Value = "43.230.08"
$ws.Range("E18").Value = "  +1.86%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.24"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.49%  "

# Row 20
$ws.Range("E20").Value = "  +6.44%  "

# Row 21
$ws.Range("E21").Value = "  +0.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.63"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.53%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
$ws.Range("E24").Value = "  -2.23%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.06%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.25%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.87"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.37%  "

# Row 28
$ws.Range("E28").Value = "  +6.76%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.15"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.26%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.62"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.47%  "

# Row 31
$ws.Range("E31").Value = "  +0.00%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.11"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.60%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0736"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.17"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.08%  "

# Row 35
$ws.Range("E35").Value = "  +7.12%  "

# Row 36
$ws.Range("E36").Value = "  +2.58%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.33"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.31"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.86%  "

# Row 39
$ws.Range("E39").Value = "  +4.70%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.43"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +12.57%  "

# Row 41
$ws.Range("E41").Value = "  +0.28%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.05"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -35.88%  "

# Row 43
$ws.Range("D43").Value = "1.958.15"
$ws.Range("E43").Value = "  +0.15%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0281"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.93%  "

# Row 45
$ws.Range("E45").Value = "  +1.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.73%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.29"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -10.17%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.91"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.31%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.52"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.06%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.04"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.85%  "

# Row 51
$ws.Range("E51").Value = "  +1.32%  "
